# Update MS data checking
# Append a new day (2020-04-08) of rows to the "Confirmados" and "Mortes"
# sheets, mirroring the layout of the existing per-day rows (date text in
# column A, per-UF integer counts in columns B:AB).

$wb = $excel.ActiveWorkbook

$newDate = "2020-04-08"

# Per-UF values (columns B..AB) in header order:
# Acre, Alagoas, Amapá, Amazonas, Bahia, Ceará, Distrito Federal,
# Espírito Santo, Goiás, Maranhão, Mato Grosso, Mato Grosso do Sul,
# Minas Gerais, Pará, Paraíba, Paraná, Pernambuco, Piauí, Rio de Janeiro,
# Rio Grande do Norte, Rio Grande do Sul, Rondônia, Roraima, Santa Catarina,
# São Paulo, Sergipe, Tocantins
$confirmados = @(54,37,107,804,497,1291,509,227,158,230,90,85,614,167,41,539,401,31,1938,261,555,18,49,457,6708,36,23)
$mortes      = @(2,2,2,30,15,43,12,6,7,11,1,2,14,6,4,17,46,5,106,11,9,1,1,15,428,4,0)

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

function Add-DayRow($sheetName, $rowNum, $values) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Write the date as literal text (matching the existing rows, which
    # store dates as shared strings rather than Excel date serials).
    # Using a formula + paste-values round trip avoids Excel's automatic
    # "looks like a date" literal-entry conversion.
    $cellRef = "A" + $rowNum
    $ws.Range($cellRef).Formula = "=""" + $newDate + """"
    $ws.Range($cellRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = 0

    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ref = $cols[$i] + $rowNum
        $ws.Range($ref).Value = $values[$i]
    }
}

Add-DayRow "Confirmados" 45 $confirmados
Add-DayRow "Mortes" 45 $mortes

Write-Output "Added 2020-04-08 row to Confirmados and Mortes sheets"
